$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Edit 1: " and " -> ". "
#   ("...does not cover the internals and structure of CMake and Ccola..."
#    becomes "...structure of CMake. Ccola...")
# ---------------------------------------------------------------
$a1 = $d.Content
$f1 = $a1.Find
$f1.Text = "structure of CMake"
[void]$f1.Execute()
$pos1 = $a1.End

$a1b = $d.Content
$f1b = $a1b.Find
$f1b.Text = "Ccola Programmer Manual."
[void]$f1b.Execute()
$paraEnd1 = $a1b.End

$sub1 = $d.Range($pos1, $paraEnd1)
$fs1 = $sub1.Find
$fs1.Text = " and "
$fs1.Replacement.Text = ". "
[void]$fs1.Execute($fs1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $fs1.Replacement.Text, 1)

# ---------------------------------------------------------------
# Edit 2: "; those are covered in the " -> " is a "
# ---------------------------------------------------------------
$a2 = $d.Content
$f2 = $a2.Find
$f2.Text = "CMake. Ccola"
[void]$f2.Execute()
$pos2 = $a2.End

$a2b = $d.Content
$f2b = $a2b.Find
$f2b.Text = "Ccola Programmer Manual."
[void]$f2b.Execute()
$paraEnd2 = $a2b.End

$sub2 = $d.Range($pos2, $paraEnd2)
$fs2 = $sub2.Find
$fs2.Text = "; those are covered in the "
$fs2.Replacement.Text = " is a "
[void]$fs2.Execute($fs2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $fs2.Replacement.Text, 1)

# ---------------------------------------------------------------
# Edit 3: second "CMake" -> "Cmake" (typo, as dictated by the target)
# ---------------------------------------------------------------
$a3 = $d.Content
$f3 = $a3.Find
$f3.Text = "Ccola is a "
[void]$f3.Execute()
$pos3 = $a3.End

$a3b = $d.Content
$f3b = $a3b.Find
$f3b.Text = "Ccola Programmer Manual."
[void]$f3b.Execute()
$paraEnd3 = $a3b.End

$sub3 = $d.Range($pos3, $paraEnd3)
$fs3 = $sub3.Find
$fs3.Text = "CMake"
$fs3.Replacement.Text = "Cmake"
[void]$fs3.Execute($fs3.Text, $false, $false, $false, $false, $false, $true, 1, $false, $fs3.Replacement.Text, 1)

# ---------------------------------------------------------------
# Edit 4: " and Ccola Programmer Manual." ->
#   " Component layer built on top of CMAKE infrastructure, described
#    further in this document."
# ---------------------------------------------------------------
$a4 = $d.Content
$f4 = $a4.Find
$f4.Text = "Ccola is a Cmake"
[void]$f4.Execute()
$pos4 = $a4.End

$a4b = $d.Content
$f4b = $a4b.Find
$f4b.Text = "Ccola Programmer Manual."
[void]$f4b.Execute()
$paraEnd4 = $a4b.End

$sub4 = $d.Range($pos4, $paraEnd4)
$fs4 = $sub4.Find
$fs4.Text = " and Ccola Programmer Manual."
$fs4.Replacement.Text = " Component layer built on top of CMAKE infrastructure, described further in this document."
[void]$fs4.Execute($fs4.Text, $false, $false, $false, $false, $false, $true, 1, $false, $fs4.Replacement.Text, 1)

# ---------------------------------------------------------------
# Edit 5: relocate the hidden "_GoBack" bookmark from the end of the
# document (after "... option:") to sit right before "further in this
# document." in the paragraph we just edited.
# ---------------------------------------------------------------
$oldBm = $d.Bookmarks("_GoBack")
$oldBm.Delete()

$a5 = $d.Content
$f5 = $a5.Find
$f5.Text = "described "
[void]$f5.Execute()
$newBmPos = $a5.End

$bmRange = $d.Range($newBmPos, $newBmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
